$d = $word.ActiveDocument

# Replace all textual occurrences of "99.95" with "99.5" throughout the document body.
$d.Content.Find.Execute("99.95", $false, $false, $false, $false, $false, $true, 1, $false, "99.5", 2)
